$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("C2").Value = 12.5
$ws.Range("C4").Value = 1.35
$ws.Range("C5").Value = 22

# Update selection to C4
$ws.Range("C4").Select()
